# Apply metadata updates to the "Metadata" sheet of the CodeSystem workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# URL: https://hl7.fr/fhir/fr/medication/CodeSystem/fr-medication-history-source-type
#   -> https://hl7.fr/ig/fhir/medication/CodeSystem/fr-medication-history-source-type
$ws.Range("B2").Value = "https://hl7.fr/ig/fhir/medication/CodeSystem/fr-medication-history-source-type"

# Name: FrMedicationHistorySourceType -> FRMedicationHistorySourceType
$ws.Range("B4").Value = "FRMedicationHistorySourceType"

# Date: 2025-04-10T15:35:36+00:00 -> 2026-01-15T08:54:26+00:00
$ws.Range("B8").Value = "2026-01-15T08:54:26+00:00"

# Jurisdiction: (empty) -> FRANCE
$ws.Range("B11").Value = "FRANCE"
